$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("PressureMean")
$ws.Range("B2").Value = 0.7113447999999999
$ws.Range("C2").Value = 0.587032
$ws.Range("B3").Value = 0.6096956486
$ws.Range("C3").Value = 0.499592544
$ws.Range("B4").Value = 0.5250332
$ws.Range("C4").Value = 0.4322551999999999
$ws.Range("B5").Value = 0.9643600000000001
$ws.Range("C5").Value = 0.90496
$ws.Range("B6").Value = 0.8904

$ws = $wb.Worksheets.Item("ActivityContributionsError")
$ws.Range("D2").Value = 0.009833094909913095
$ws.Range("D3").Value = 0.01503795225662456
$ws.Range("D4").Value = 0.008203096626521822
$ws.Range("D5").Value = 0.01231084081331761
$ws.Range("D6").Value = 0.01217620155530652
$ws.Range("D7").Value = 0.0004063846268960041
$ws.Range("D8").Value = 0.001230300045353961
$ws.Range("D9").Value = 0.004203198890535038
$ws.Range("D10").Value = 0.005436662274549129
$ws.Range("D11").Value = 0.01532015968747565
$ws.Range("D12").Value = 0.01390568315370264

$ws = $wb.Worksheets.Item("PressureContributionsMean")
$ws.Range("D2").Value = 0.6727549802662831
$ws.Range("D3").Value = 0.2113769774816732
$ws.Range("D4").Value = 0.1158680422520436
$ws.Range("D5").Value = 0.83236812445981
$ws.Range("D6").Value = 0.1676318755401902
$ws.Range("D7").Value = 0.6343263002050215
$ws.Range("D8").Value = 0.250559293105184
$ws.Range("D9").Value = 0.1151144066897944
$ws.Range("D10").Value = 0.8352891734856789
$ws.Range("D11").Value = 0.1647108265143212

$ws = $wb.Worksheets.Item("PressureContributionsError")
$ws.Range("D2").Value = 0.009322458761353375
$ws.Range("D3").Value = 0.009657071528972789
$ws.Range("D4").Value = 0.001632923101321011
$ws.Range("D5").Value = 0.002487358973910752
$ws.Range("D6").Value = 0.002487358973910721
$ws.Range("D7").Value = 0.008045319065386506
$ws.Range("D8").Value = 0.006492128472678519
$ws.Range("D9").Value = 0.001944735896859549
$ws.Range("D10").Value = 0.004606348893236272
$ws.Range("D11").Value = 0.004606348893236218

$ws = $wb.Worksheets.Item("PressureError")
$ws.Range("B2").Value = 0.01215713906410551
$ws.Range("C2").Value = 0.01224750253725224
$ws.Range("B3").Value = 0.009344207559866452
$ws.Range("C3").Value = 0.01188433385159427
$ws.Range("B4").Value = 0.01444813800214024
$ws.Range("C4").Value = 0.008498693898084178
$ws.Range("B5").Value = 0.0005400000000000001
$ws.Range("C5").Value = 0.001439999999999995
$ws.Range("B6").Value = 0.01626666666666667

$ws = $wb.Worksheets.Item("TPLMean")
$ws.Range("B2").Value = 0.6131565049253417
$ws.Range("C2").Value = 0.503246469473797
$ws.Range("B3").Value = 0.6223306234372835
$ws.Range("C3").Value = 0.5057565134872066

$ws = $wb.Worksheets.Item("TPLError")
$ws.Range("B2").Value = 0.008708835448271574
$ws.Range("C2").Value = 0.007163802701522741
$ws.Range("B3").Value = 0.00923368994353385
$ws.Range("C3").Value = 0.01066197055715727

$ws = $wb.Worksheets.Item("MeasureEffectsMean")
$ws.Range("E2").Value = 0.411
$ws.Range("E3").Value = 0.5870000000000001
$ws.Range("E4").Value = 0.354
$ws.Range("E6").Value = 0.344
$ws.Range("E7").Value = 0.274
$ws.Range("E8").Value = 0.396
$ws.Range("E9").Value = 0.424
$ws.Range("E10").Value = 0.599
$ws.Range("E11").Value = 0.659

$ws = $wb.Worksheets.Item("MeasureEffectsError")
$ws.Range("E2").Value = 0.0267685553505518
$ws.Range("E3").Value = 0.03732886878066954
$ws.Range("E4").Value = 0.01967513941783161
$ws.Range("E6").Value = 0.04292629341868066
$ws.Range("E7").Value = 0.04066666666666666
$ws.Range("E8").Value = 0.005999999999999998
$ws.Range("E9").Value = 0.01351542328847553
$ws.Range("E10").Value = 0.02639233895575676
$ws.Range("E11").Value = 0.01940503943710385

$ws = $wb.Worksheets.Item("ActivityContributionsMean")
$ws.Range("D2").Value = 0.4361746830281742
$ws.Range("D3").Value = 0.3157890423588068
$ws.Range("D4").Value = 0.5381225776204775
$ws.Range("D5").Value = 0.4987906758790127
$ws.Range("D6").Value = 0.3241463999621798
$ws.Range("D7").Value = 0.2741279449337248
$ws.Range("D8").Value = 0.2264675296970426
$ws.Range("D9").Value = 0.2152965150799589
$ws.Range("D10").Value = 0.1993508173286767
$ws.Range("D11").Value = 0.4248273930347688
$ws.Range("D12").Value = 0.303509595157636
